$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly data rows appended below the existing table (rows 124-126).
$rows = @(
    @{ Row = 124; A = "Eric";   C = "Workout"; D = 60; E = 0; F = 0; G = 13; H = 39; I = 7;  J = 2; K = 0; L = "Sauntering Hippo"; M = 3 },
    @{ Row = 125; A = "Steven"; C = "Workout"; D = 26; E = 0; F = 0; G = 11; H = 13; I = 3;  J = 0; K = 0; L = "Wily Hyena";       M = 3 },
    @{ Row = 126; A = "Phil";   C = "Workout"; D = 65; E = 0; F = 0; G = 12; H = 36; I = 16; J = 1; K = 0; L = "Sauntering Hippo"; M = 3 }
)

foreach ($r in $rows) {
    $n = $r.Row

    $ws.Cells.Item($n, 1).Value = $r.A

    # Column B holds the workout date; copy the date cell directly above so
    # the new cell inherits the existing short-date style (matches the
    # m/dd/yy numeric formatting already used for every other date cell).
    $ws.Cells.Item($n - 1, 2).Copy($ws.Cells.Item($n, 2))
    $ws.Cells.Item($n, 2).Value = 45471

    $ws.Cells.Item($n, 3).Value = $r.C
    $ws.Cells.Item($n, 4).Value = $r.D
    $ws.Cells.Item($n, 5).Value = $r.E
    $ws.Cells.Item($n, 6).Value = $r.F
    $ws.Cells.Item($n, 7).Value = $r.G
    $ws.Cells.Item($n, 8).Value = $r.H
    $ws.Cells.Item($n, 9).Value = $r.I
    $ws.Cells.Item($n, 10).Value = $r.J
    $ws.Cells.Item($n, 11).Value = $r.K
    $ws.Cells.Item($n, 12).Value = $r.L
    $ws.Cells.Item($n, 13).Value = $r.M
}

# Keep the active selection in sync with the newly-added last row, mirroring
# where Excel leaves the cursor after entering this data.
$ws.Range("A127").Select()
